# Applies the "ok 78.5, add interact file" edit:
#  - Adds 14 new rows (231-244) to Sheet1 that repeat the "append chap them
#    noi them / ..." phrase set from row 229, except columns Q and S which
#    reference two brand new vocabulary entries:
#       Q -> "affect anh huong den"
#       S -> "difinitely chac chan"
#  - Moves the sheet view's scroll position / active selection down to
#    around row 225-244 (topLeftCell P225, active cell U244).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Same 20 column values as row 229, but with Q and S pointing at the two
# newly introduced words.
$rowValues = @(
    "append chap them noi them",
    "instance truong hop vi du",
    "emit xong len phat ra",
    "enroll ghi danh",
    "current hien hanh",
    "elegant thanh lich",
    "simultaneously dong thoi",
    "critical phan doan phe binh chi trich",
    "distinct khac biet",
    "protocol giao thuc",
    "sequel phan tiep theo",
    "bind troi buoc",
    "identified xac dinh",
    "propagate lan ra truyen ra",
    "handled xu ly",
    "nowadays ngay nay",
    "affect anh huong den",
    "explicitly ro rang",
    "difinitely chac chan",
    "simulate mo phong"
)

for ($r = 231; $r -le 244; $r++) {
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $rowValues[$i]
    }
}

# Scroll the view down toward the newly added rows and move the selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 225
$win.ScrollColumn = 16
$ws.Range("U244").Select()
